# Update per-player match stats on Sheet1 (LKS.xlsx roster/stats sheet).
# Values are stored as inline-string (text) cells, so we assign strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D. Arndt (row 2)
$ws.Range("E2").Value = "540"
$ws.Range("F2").Value = "6"
$ws.Range("G2").Value = "6"

# A. Bobek (row 3)
$ws.Range("J3").Value = "1"

# M. Bakowicz (row 6)
$ws.Range("J6").Value = "1"

# M. Dabrowski (row 7)
$ws.Range("E7").Value = "630"
$ws.Range("F7").Value = "7"
$ws.Range("G7").Value = "7"

# K. Dankowski (row 8)
$ws.Range("E8").Value = "383"
$ws.Range("F8").Value = "7"
$ws.Range("G8").Value = "3"

# O. Koprowski (row 9)
$ws.Range("E9").Value = "104"
$ws.Range("F9").Value = "2"
$ws.Range("G9").Value = "1"

# M. Lorenc (row 10)
$ws.Range("E10").Value = "91"
$ws.Range("F10").Value = "3"
$ws.Range("G10").Value = "1"
$ws.Range("I10").Value = "1"
$ws.Range("K10").Value = "1"

# Nacho Monsalve (row 12)
$ws.Range("E12").Value = "630"
$ws.Range("F12").Value = "7"
$ws.Range("G12").Value = "7"

# M. Wszolek (row 14)
$ws.Range("J14").Value = "3"

# B. Biel (row 15)
$ws.Range("E15").Value = "434"
$ws.Range("F15").Value = "7"
$ws.Range("G15").Value = "6"
$ws.Range("I15").Value = "5"

# K. Ibe-Torti (row 17)
$ws.Range("E17").Value = "225"
$ws.Range("F17").Value = "7"
$ws.Range("H17").Value = "6"
$ws.Range("J17").Value = "6"

# D. Kort (row 18)
$ws.Range("E18").Value = "400"
$ws.Range("F18").Value = "6"
$ws.Range("G18").Value = "5"
$ws.Range("I18").Value = "5"

# M. Kowalczyk (row 19)
$ws.Range("E19").Value = "422"
$ws.Range("F19").Value = "7"
$ws.Range("G19").Value = "5"
$ws.Range("I19").Value = "3"
$ws.Range("L19").Value = "4"

# J. Kuzma (row 20)
$ws.Range("E20").Value = "223"
$ws.Range("F20").Value = "6"
$ws.Range("H20").Value = "4"
$ws.Range("J20").Value = "5"

# Javi Moreno (row 21)
$ws.Range("J21").Value = "7"

# V. Okhronchuk (row 23)
$ws.Range("E23").Value = "167"
$ws.Range("F23").Value = "5"
$ws.Range("H23").Value = "3"
$ws.Range("J23").Value = "5"

# Pirulo (row 25)
$ws.Range("E25").Value = "582"
$ws.Range("F25").Value = "7"
$ws.Range("G25").Value = "7"
$ws.Range("I25").Value = "3"
$ws.Range("K25").Value = "4"

# M. Trabka (row 26)
$ws.Range("E26").Value = "542"
$ws.Range("F26").Value = "7"
$ws.Range("H26").Value = "1"
$ws.Range("J26").Value = "1"

# N. Balongo (row 27)
$ws.Range("E27").Value = "436"
$ws.Range("F27").Value = "7"
$ws.Range("G27").Value = "5"

# P. Janczukowicz (row 28)
$ws.Range("E28").Value = "191"
$ws.Range("F28").Value = "5"
$ws.Range("H28").Value = "3"
$ws.Range("J28").Value = "4"
